$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": update the notes explaining the dispatch-priority rules.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A7").Value = "We assign priority 1 to natural gas peaker and petroleum-fired plants, which are the"
$about.Range("A8").Value = "only types for which a non-zero quantity is specified for guaranteed dispatch in the BAU case"
$about.Range("A9").Value = "for the United States.  We arbitrarily assign priority 2 to all other plant types."

# ---------------------------------------------------------------------------
# Sheet "BDPbES": update priorities for natural gas peaker / petroleum to 1,
# and append new plant types (rows 18-24) with priority 2.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BDPbES")

# petroleum (row 11) and natural gas peaker (row 12) now get priority 1
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1

# New plant type rows, each with priority 2 filled across all year columns
$newTypes = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen combustion turbine",
    "hydrogen combined cycle"
)

$row = 18
foreach ($name in $newTypes) {
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = 2
    $ws.Range("C$row`:AK$row").Formula = '=$B' + $row
    $row++
}

# The last two new rows (hydrogen combustion turbine / hydrogen combined
# cycle) get a distinct label style: explicit black font + vertically
# centered alignment. Apply it to A23 first, then replicate the exact same
# format onto A24 via a formats-only paste (keeps a single new style entry
# instead of one per cell).
$ws.Range("A23").Font.Color = 0
$ws.Range("A23").VerticalAlignment = -4108

$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
